$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 191: 04-10-2021
$ws.Range("A191").NumberFormat = "@"
$ws.Range("A191").Value = "04-10-2021"
$ws.Range("A191").ClearFormats()
$ws.Range("B191").Value = 116.83
$ws.Range("C191").Value = 190.46
$ws.Range("D191").Value = 103.06

# Row 192: 05-10-2021
$ws.Range("A192").NumberFormat = "@"
$ws.Range("A192").Value = "05-10-2021"
$ws.Range("A192").ClearFormats()
$ws.Range("B192").Value = 116.94
$ws.Range("C192").Value = 191.12
$ws.Range("D192").Value = 103.12
